$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71

# Text columns (Date, Time, Weekday, Week) - force text so Excel does not
# auto-convert date/time-looking or numeric-looking strings.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-15"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "12:25:47"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "06"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns (Beijing .. Wuhan)
$ws.Cells.Item($row, 5).Value = 126870
$ws.Cells.Item($row, 6).Value = 139877
$ws.Cells.Item($row, 7).Value = 170532
$ws.Cells.Item($row, 8).Value = 159538
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145048
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192220
$ws.Cells.Item($row, 14).Value = 115128
$ws.Cells.Item($row, 15).Value = 45009
$ws.Cells.Item($row, 16).Value = 28740
$ws.Cells.Item($row, 17).Value = 65883
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 45962
$ws.Cells.Item($row, 20).Value = -1
